$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-09-06 Saturday" "2025-09-07 Sunday"

Replace-Text "306×9=" "761×7="
Replace-Text "699×2=" "684×4="
Replace-Text "499×3=" "164×8="
Replace-Text "655×4=" "665×8="
Replace-Text "879×6=" "197×9="

Replace-Text "884×4=" "981×4="
Replace-Text "135×9=" "165×8="
Replace-Text "118×3=" "874×3="
Replace-Text "935×4=" "747×2="
Replace-Text "656×5=" "391×9="

Replace-Text "166×8=" "114×6="
Replace-Text "545×4=" "774×3="
Replace-Text "587×3=" "371×7="
Replace-Text "346×5=" "134×3="
Replace-Text "750×6=" "761×7="

Replace-Text "813×7=" "869×2="
Replace-Text "859×8=" "539×9="
Replace-Text "334×9=" "416×7="
Replace-Text "907×2=" "151×4="
Replace-Text "492×9=" "856×6="

Replace-Text "860×5=" "436×5="
Replace-Text "951×5=" "509×3="
Replace-Text "206×7=" "704×5="
Replace-Text "965×6=" "986×9="
Replace-Text "886×9=" "164×5="
